# Update the cryptos price/volume table (columns D and E) with refreshed
# values, as produced by the scheduled GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.666.36'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '3.589.79'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '609.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('D7').Value = '3.588.12'
$ws.Range('E7').Value = '  +0.63%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.491'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.29%  '
$ws.Range('E10').Value = '  +0.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.99'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('E12').Value = '  +1.17%  '
$ws.Range('D13').Value = '4.199.82'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '30.14'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').Value = '3.586.54'
$ws.Range('E16').Value = '  +0.55%  '
$ws.Range('D17').Value = '66.748.04'
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.47'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.32'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '431.84'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.626'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.76%  '
$ws.Range('D25').Value = '3.734.69'
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  -0.89%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('E30').Value = '  -0.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').Value = '3.586.53'
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '25.53'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.38%  '
$ws.Range('E34').Value = '  -2.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.154'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.09%  '
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('E39').Value = '  -0.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '173.53'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0858'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.25'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.897'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.92'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('E45').Value = '  +0.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.56'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.91%  '
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.99'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.23'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.12%  '
